# Refresh market-price snapshot values (currentAveragePrice.. / LevePrice.. / LeveProfit..,
# columns H:N) for the affected Leve rows across all crafting-class sheets.
# Source data only - no formulas live in these columns, so cells are written directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13: The Hexster Runoff / Maple Picatrix
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 2259858
$ws.Range("J17").Value = 2292140
$ws.Range("L17").Value = 6876420
$ws.Range("N17").Value = -6876756

# Row 46: Always Have an Exit Plan / Poisoning Potion
$ws.Range("H46").Value = 1321.5714
$ws.Range("I46").Value = 583.6667
$ws.Range("J46").Value = 1875
$ws.Range("K46").Value = 1751.0001
$ws.Range("L46").Value = 5625
$ws.Range("M46").Value = -1632.0001
$ws.Range("N46").Value = -5863

# Row 60: Make Up Your Mind or Else / Potent Poisoning Potion
$ws.Range("H60").Value = 1321.5714
$ws.Range("I60").Value = 583.6667
$ws.Range("J60").Value = 1875
$ws.Range("K60").Value = 1751.0001
$ws.Range("L60").Value = 5625
$ws.Range("M60").Value = -1267.0001
$ws.Range("N60").Value = -6593

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 23812042
$ws.Range("I132").Value = 2404.125
$ws.Range("J132").Value = 55558224
$ws.Range("K132").Value = 7212.375
$ws.Range("L132").Value = 166674672
$ws.Range("M132").Value = -4682.375
$ws.Range("N132").Value = -166679732

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 6784.55
$ws.Range("I61").Value = 6784.55
$ws.Range("K61").Value = 6784.55
$ws.Range("M61").Value = -6572.55

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 17858638
$ws.Range("I74").Value = 1277.1904
$ws.Range("J74").Value = 71430720
$ws.Range("K74").Value = 1277.1904
$ws.Range("L74").Value = 71430720
$ws.Range("M74").Value = -403.1904
$ws.Range("N74").Value = -71432468

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 17858638
$ws.Range("I77").Value = 1277.1904
$ws.Range("J77").Value = 71430720
$ws.Range("K77").Value = 6385.951999999999
$ws.Range("L77").Value = 357153600
$ws.Range("M77").Value = -2017.951999999999
$ws.Range("N77").Value = -357162336

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1951284.2
$ws.Range("I102").Value = 2059522.2
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2059522.2
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -2057900.2
$ws.Range("N102").Value = -6244

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2880.1462
$ws.Range("I132").Value = 2072.125
$ws.Range("J132").Value = 5753.1113
$ws.Range("K132").Value = 6216.375
$ws.Range("L132").Value = 17259.3339
$ws.Range("M132").Value = -3686.375
$ws.Range("N132").Value = -22319.3339

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 6784.55
$ws.Range("I136").Value = 6784.55
$ws.Range("K136").Value = 20353.65
$ws.Range("M136").Value = -17803.65

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 3795.353
$ws.Range("I134").Value = 4189.564
$ws.Range("J134").Value = 2514.1667
$ws.Range("K134").Value = 12568.692
$ws.Range("L134").Value = 7542.500100000001
$ws.Range("M134").Value = -10033.692
$ws.Range("N134").Value = -12612.5001

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 6527613.5
$ws.Range("I31").Value = 1707.4419
$ws.Range("J31").Value = 14322446
$ws.Range("K31").Value = 1707.4419
$ws.Range("L31").Value = 14322446
$ws.Range("M31").Value = -1412.4419
$ws.Range("N31").Value = -14323036

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 6527613.5
$ws.Range("I34").Value = 1707.4419
$ws.Range("J34").Value = 14322446
$ws.Range("K34").Value = 1707.4419
$ws.Range("L34").Value = 14322446
$ws.Range("M34").Value = -1505.4419
$ws.Range("N34").Value = -14322850

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2874661.2
$ws.Range("I58").Value = 3968845.8
$ws.Range("J58").Value = 2427.3125
$ws.Range("K58").Value = 3968845.8
$ws.Range("L58").Value = 2427.3125
$ws.Range("M58").Value = -3968642.8
$ws.Range("N58").Value = -2833.3125

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 2527.1333
$ws.Range("I122").Value = 2072
$ws.Range("J122").Value = 3778.75
$ws.Range("K122").Value = 6216
$ws.Range("L122").Value = 11336.25
$ws.Range("M122").Value = -3766
$ws.Range("N122").Value = -16236.25

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3334892
$ws.Range("I132").Value = 4349204.5
$ws.Range("J132").Value = 2151.4285
$ws.Range("K132").Value = 13047613.5
$ws.Range("L132").Value = 6454.2855
$ws.Range("M132").Value = -13045083.5
$ws.Range("N132").Value = -11514.2855

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 7409543.5
$ws.Range("I134").Value = 13891855
$ws.Range("J134").Value = 1188
$ws.Range("K134").Value = 41675565
$ws.Range("L134").Value = 3564
$ws.Range("M134").Value = -41673030
$ws.Range("N134").Value = -8634

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2874661.2
$ws.Range("I136").Value = 3968845.8
$ws.Range("J136").Value = 2427.3125
$ws.Range("K136").Value = 11906537.4
$ws.Range("L136").Value = 7281.9375
$ws.Range("M136").Value = -11903987.4
$ws.Range("N136").Value = -12381.9375

$ws = $wb.Worksheets.Item("CUL")
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 2084098.9
$ws.Range("I131").Value = 3704232
$ws.Range("J131").Value = 1070.8096
$ws.Range("K131").Value = 11112696
$ws.Range("L131").Value = 3212.4288
$ws.Range("M131").Value = -11107656
$ws.Range("N131").Value = -13292.4288

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 21787676
$ws.Range("J132").Value = 32681114
$ws.Range("L132").Value = 294130026
$ws.Range("N132").Value = -294135086

# Row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 50524.39
$ws.Range("I133").Value = 95823.73
$ws.Range("J133").Value = 9000
$ws.Range("K133").Value = 287471.19
$ws.Range("L133").Value = 27000
$ws.Range("M133").Value = -282411.19
$ws.Range("N133").Value = -37120

$ws = $wb.Worksheets.Item("GSM")
# Row 14: All That Glitters / Copper Ear Cuffs
$ws.Range("H14").Value = 3672000
$ws.Range("I14").Value = 6601600
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 6601600
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -6601432
$ws.Range("N14").Value = -10336

$ws = $wb.Worksheets.Item("LTW")
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 299484.34
$ws.Range("I82").Value = 501374.16
$ws.Range("K82").Value = 501374.16
$ws.Range("M82").Value = -501013.16

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 299484.34
$ws.Range("I85").Value = 501374.16
$ws.Range("K85").Value = 501374.16
$ws.Range("M85").Value = -500126.16

# Row 104: Brace Yourselves / Gazelleskin Bracers of Fending
$ws.Range("H104").Value = 9000
$ws.Range("J104").Value = 9000
$ws.Range("L104").Value = 9000
$ws.Range("N104").Value = -15988

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 11351386
$ws.Range("I132").Value = 13679321
$ws.Range("J132").Value = 2705
$ws.Range("K132").Value = 41037963
$ws.Range("L132").Value = 8115
$ws.Range("M132").Value = -41035433
$ws.Range("N132").Value = -13175

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 37043304
$ws.Range("I136").Value = 8937.454
$ws.Range("J136").Value = 95240170
$ws.Range("K136").Value = 26812.362
$ws.Range("L136").Value = 285720510
$ws.Range("M136").Value = -24262.362
$ws.Range("N136").Value = -285725610
